$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.296.31"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "2.353.63"
$ws.Range("E3").Value = "  +6.26%  "
$ws.Range("E4").Value = "  -0.30%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "313.62"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +5.90%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "109.69"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.66%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.644"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.29%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  +7.53%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "43.11"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0939"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.24%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.86"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.15%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.05"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +9.39%  "
$ws.Range("E14").Value = "  +2.17%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "16.43"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +9.86%  "
$ws.Range("D16").Value = "2.704.90"
$ws.Range("E16").Value = "  +6.23%  "
$ws.Range("D17").Value = "2.475.85"
$ws.Range("E17").Value = "  +10.90%  "
$ws.Range("D18").Value = "43.279.87"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("E20").Value = "  -1.50%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "75.49"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +4.27%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.45"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +12.19%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "255.13"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +11.91%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.13"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.71%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "12.06"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.14%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "39.27"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("E29").Value = "  +1.12%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "22.37"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +7.07%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "174.02"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  -0.56%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0931"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.69%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.03"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +9.31%  "
$ws.Range("E35").Value = "  +5.82%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.99"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.15"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0377"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.39%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.70%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +11.74%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "72.93"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.64%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.49"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +14.83%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.234"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "12.79"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "5.63"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.36%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.31"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +11.25%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "111.27"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +7.79%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.463"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +6.73%  "
